$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.323.09'
$ws.Range("E2").Value = '  +0.54%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.667.54'
$ws.Range("E3").Value = '  +0.78%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.008'
$ws.Range("E4").Value = '  +0.23%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '220.65'
$ws.Range("E5").Value = '  +1.35%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5311'
$ws.Range("E6").Value = '  +0.01%  '

# Row 7
$ws.Range("E7").Value = '  +0.22%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2649'
$ws.Range("E8").Value = '  +0.87%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06373'
$ws.Range("E9").Value = '  +0.55%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.86'
$ws.Range("E10").Value = '  +2.29%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07830'
$ws.Range("E11").Value = '  +0.30%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.520'
$ws.Range("E12").Value = '  -0.02%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.661.35'
$ws.Range("E13").Value = '  +2.13%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.895.91'
$ws.Range("E14").Value = '  +0.77%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5600'
$ws.Range("E15").Value = '  +1.97%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₅8168'
$ws.Range("E16").Value = '  +0.24%  '

# Row 17
$ws.Range("E17").Value = '  +0.70%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '26.339.36'
$ws.Range("E18").Value = '  +0.78%  '

# Row 19
$ws.Range("E19").Value = '  +0.25%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.721'
$ws.Range("E20").Value = '  +2.84%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '197.26'
$ws.Range("E21").Value = '  +3.22%  '

# Row 22
$ws.Range("E22").Value = '  +1.76%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.047'
$ws.Range("E23").Value = '  +0.53%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.009'
$ws.Range("E24").Value = '  +0.17%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.24'
$ws.Range("E25").Value = '  +0.70%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1220'
$ws.Range("E26").Value = '  +0.28%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.249'
$ws.Range("E27").Value = '  +0.79%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '16.16'

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.507'
$ws.Range("E29").Value = '  +2.28%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05888'
$ws.Range("E30").Value = '  +2.44%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.283'
$ws.Range("E31").Value = '  +0.76%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.546'
$ws.Range("E32").Value = '  -0.17%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.330'
$ws.Range("E33").Value = '  +2.01%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.604'
$ws.Range("E34").Value = '  +0.98%  '

# Row 35
$ws.Range("E35").Value = '  +0.73%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9609'
$ws.Range("E36").Value = '  +1.16%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.431'
$ws.Range("E37").Value = '  +0.35%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5821'
$ws.Range("E38").Value = '  +1.26%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01614'
$ws.Range("E39").Value = '  +0.81%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.957'
$ws.Range("E40").Value = '  +2.89%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.078.61'
$ws.Range("E41").Value = '  +3.34%  '

# Row 42
$ws.Range("E42").Value = '  +0.78%  '

# Row 43
$ws.Range("E43").Value = '  +0.26%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '102.83'

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.805.90'

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '58.45'
$ws.Range("E46").Value = '  +3.12%  '

# Row 47
$ws.Range("B47").Value = 'Frax'
$ws.Range("C47").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.016'
$ws.Range("E47").Value = '  +1.39%  '

# Row 48
$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0₈105'
$ws.Range("E48").Value = '  +0.44%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4407'
$ws.Range("E49").Value = '  +1.12%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.060'
$ws.Range("E50").Value = '  +2.50%  '
